$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B13").Value = "6b5bef8bc5d42f9d30b95d561587e88f"
$ws.Range("B17").Value = "09088a9949ad1a8f49e48e974ac83087"
$ws.Range("B26").Value = "27ea83e59e475359c7067824c461be37"
$ws.Range("B36").Value = "cf0dc0e7bccb941d33745a261424c664"
$ws.Range("B56").Value = "ad89823b30623527e3baef2253bbb1b1"
$ws.Range("B82").Value = "c4ba01d46b527d0a4e04edce448dd616"
$ws.Range("B123").Value = "25fbfe6bc6a3bd0a2084d71b2d872320"
$ws.Range("B126").Value = "4b993743580714c27f04099bacdaa06f"
$ws.Range("B136").Value = "fbec840e02f535bdf0b4b1cd69c2a46a"
$ws.Range("B166").Value = "2b0269c464b59fe3dd3489ae566592a3"
$ws.Range("B185").Value = "c861d4fd025db76c7afc46ce10a49f3f"
$ws.Range("B196").Value = "0fa317cd20a5b49aad28226bf6df9533"
$ws.Range("B203").Value = "bd9f44647569b5e2f4c5ad711528b41a"
$ws.Range("B218").Value = "bc91d4b0be812322444314f5443b41bc"
$ws.Range("B341").Value = "9957746ca62d480b301fca0d3984ec02"
$ws.Range("B347").Value = "eb4717392d96cb5d186e478cf51f188e"
$ws.Range("B433").Value = "0841f66eec1f7caf51680bed6f5054c6"
$ws.Range("B498").Value = "830064ea44c0d05c93a59723f7118ea2"
$ws.Range("B512").Value = "d1484942ef08d715b335040b3a5e7c24"
$ws.Range("B527").Value = "067ca15bdc2bb75fdda992d4b749a669"
$ws.Range("B535").Value = "6c1aa9bb02a7b3c0cdde7e4b3fc07dbf"
$ws.Range("B536").Value = "e36cfc8d647d2ba777e89889eb5fd238"
$ws.Range("B546").Value = "d798c3c15221dcec831d881a939029c2"
$ws.Range("B566").Value = "282e369444428aa780db634d3867d417"
$ws.Range("B597").Value = "b7ae7df38ad9ddc0f435e255584d2b1b"
$ws.Range("B598").Value = "819a318a42b307090b15e32fe333138c"
$ws.Range("B607").Value = "ba39438bc23efc7a39321ed0bae1a377"
$ws.Range("B659").Value = "6da39242dc3e342481be3b884dfb17d8"
$ws.Range("B690").Value = "79ec90319fcd89f911fb8849489458d2"
$ws.Range("B729").Value = "69e752efebe41c5f6e6b78c9585d5f64"
$ws.Range("B751").Value = "952382d47d10bd78b5a7cebabbc8493e"
$ws.Range("B754").Value = "00befa18721ee13a08f8b9160e59350c"
$ws.Range("B758").Value = "c9d8ad33e3779fce70cd4c812e3a58e9"
$ws.Range("B767").Value = "9c12e2aaa7853b75756fbbc3d8715c29"
$ws.Range("B879").Value = "43057fba58fa2f6a9cc65bd8ce502873"
$ws.Range("B883").Value = "8f893dcb52b1e47d448c2b51b718df3e"
